# Generate Report for Handoff
# This script updates the localization-status workbook so that the two
# tracked files (previously handed back) are now shown as newly queued for
# handoff: new source file identifiers, status "Ready for handoff", fresh
# handoff-only file/datetime info, and the old handback-specific columns
# (Latest Target File / Latest Handback File) cleared out.

$wb = $excel.ActiveWorkbook

$oldFile1 = "3025d094-15e2-455b-a109-7a5d624e4eb5.md"
$oldFile2 = "7d15c9a5-b4a7-4b0c-9e48-62ddeb6b66ff.md"

$newFile1 = "8bf89a9d-2aaf-449d-921c-57af386ba3c2.md"
$newFile2 = "ffff141e73cc-5ebe-4266-951d-c32a64d51a45.md"

$newStatus = "Ready for handoff"

$newXlfZh = "8bf89a9d-2aaf-449d-921c-57af386ba3c2.41728e5077dc227b602135e289cdca0fded0e11e.zh-cn.xlf"
$newXlfDe = "8bf89a9d-2aaf-449d-921c-57af386ba3c2.41728e5077dc227b602135e289cdca0fded0e11e.de-de.xlf"

$newDateZh = "2016-03-08 14:55:18"
$newDateDe = "2016-03-08 14:55:23"
$epochDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/e2e/" + $newFile1, "", "", $newFile1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/e2e/" + $newFile2, "", "", $newFile2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/.localization-config", "", "", ".localization-config")

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

Write-Host "Updated Overview sheet"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

# Clear the columns that are no longer populated for the new rows
# (Latest Target File / Latest Handback File / Dependency From)
$wsZh.Range("E2").Clear()
$wsZh.Range("F2").Clear()
$wsZh.Range("I2").Clear()
$wsZh.Range("E3").Clear()
$wsZh.Range("F3").Clear()
$wsZh.Range("I3").Clear()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/e2e/" + $newFile1, "", "", $newFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56c98c74da86e600422e560e2a5e1e75ac6f14bf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $newXlfZh, "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/e2e/" + $newFile2, "", "", $newFile2)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/56c98c74da86e600422e560e2a5e1e75ac6f14bf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $newXlfZh, "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/.localization-config", "", "", ".localization-config")

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("D2").Value = $newDateZh
$wsZh.Range("G2").Value = $epochDate
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("D3").Value = $newDateZh
$wsZh.Range("G3").Value = $epochDate
$wsZh.Range("H3").Value = "Include"

Write-Host "Updated zh-cn sheet"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

# Clear the columns that are no longer populated for the new rows
# (Latest Target File / Latest Handback File / Dependency From)
$wsDe.Range("E2").Clear()
$wsDe.Range("F2").Clear()
$wsDe.Range("I2").Clear()
$wsDe.Range("E3").Clear()
$wsDe.Range("F3").Clear()
$wsDe.Range("I3").Clear()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/e2e/" + $newFile1, "", "", $newFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eac8c8fcf57cbbe5d3e3a487ea661ba77ca4be79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $newXlfDe, "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/e2e/" + $newFile2, "", "", $newFile2)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eac8c8fcf57cbbe5d3e3a487ea661ba77ca4be79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $newXlfDe, "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d5071e80a525a54540d6eb768e47d537a98e06ea/.localization-config", "", "", ".localization-config")

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("D2").Value = $newDateDe
$wsDe.Range("G2").Value = $epochDate
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("D3").Value = $newDateDe
$wsDe.Range("G3").Value = $epochDate
$wsDe.Range("H3").Value = "Include"

Write-Host "Updated de-de sheet"
